$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Price" (D) and "Volume(1h)" (E) columns of the cryptos table
# with freshly scraped figures. Price strings such as "65.568.69" use dots as
# thousands separators (not valid Excel numbers) and must stay text; a few
# updated prices (e.g. "602.20", "7.89") WOULD parse as plain numbers, so for
# those we pin the cell to text format first to keep the literal digits/zeros
# exactly as scraped.

$ws.Range("D2").Value = "65.568.69"
$ws.Range("E2").Value = "  -1.50%  "
$ws.Range("D3").Value = "3.526.00"
$ws.Range("E3").Value = "  -1.63%  "
$ws.Range("E4").Value = "  -0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.20"
$ws.Range("E5").Value = "  -1.07%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "143.75"
$ws.Range("E6").Value = "  -1.67%  "
$ws.Range("D7").Value = "3.523.55"
$ws.Range("E7").Value = "  -1.62%  "
$ws.Range("E8").Value = "  -0.19%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.516"
$ws.Range("E9").Value = "  +5.27%  "
$ws.Range("E10").Value = "  -1.99%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "7.89"
$ws.Range("E11").Value = "  -1.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.405"
$ws.Range("E12").Value = "  -2.54%  "
$ws.Range("D13").Value = "4.115.79"
$ws.Range("E13").Value = "  -1.99%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000197"
$ws.Range("E14").Value = "  -5.38%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "28.45"
$ws.Range("E15").Value = "  -5.11%  "
$ws.Range("D16").Value = "3.523.26"
$ws.Range("E16").Value = "  -1.89%  "
$ws.Range("E17").Value = "  +1.56%  "
$ws.Range("D18").Value = "65.506.48"
$ws.Range("E18").Value = "  -1.71%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.04"
$ws.Range("E19").Value = "  -3.67%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "6.20"
$ws.Range("E20").Value = "  -1.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.34"
$ws.Range("E21").Value = "  -4.29%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "416.82"
$ws.Range("E22").Value = "  -3.55%  "
$ws.Range("E23").Value = "  -3.83%  "
$ws.Range("E24").Value = "  -1.88%  "
$ws.Range("D25").Value = "3.658.21"
$ws.Range("E25").Value = "  -1.98%  "
$ws.Range("E26").Value = "  +0.09%  "
$ws.Range("E27").Value = "  -5.30%  "
$ws.Range("E28").Value = "  -2.63%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "7.78"
$ws.Range("E29").Value = "  -3.70%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.88"
$ws.Range("E30").Value = "  -4.43%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.998"
$ws.Range("E31").Value = "  -0.17%  "
$ws.Range("D32").Value = "3.526.26"
$ws.Range("E32").Value = "  -1.67%  "
$ws.Range("E33").Value = "  -1.72%  "
$ws.Range("E34").Value = "  -4.49%  "
$ws.Range("E35").Value = "  -0.01%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.31"
$ws.Range("E36").Value = "  -9.21%  "
$ws.Range("E37").Value = "  -4.31%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "174.60"
$ws.Range("E38").Value = "  +0.38%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.27"
$ws.Range("E39").Value = "  -6.30%  "
$ws.Range("E40").Value = "  -8.33%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0814"
$ws.Range("E41").Value = "  -4.51%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.06"
$ws.Range("E42").Value = "  -2.96%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.855"
$ws.Range("E43").Value = "  -4.28%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "45.05"
$ws.Range("E44").Value = "  -1.78%  "
$ws.Range("E45").Value = "  -7.94%  "
$ws.Range("E46").Value = "  -0.05%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.35"
$ws.Range("E47").Value = "  -6.83%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "23.35"
$ws.Range("E48").Value = "  -1.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.01"
$ws.Range("E49").Value = "  -2.53%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.10"
$ws.Range("E50").Value = "  -7.96%  "
$ws.Range("E51").Value = "  -4.36%  "
